# Weekly update: a new price record is added for the most recent date,
# pushing the existing historical rows (17-23) down by one row (to 18-24).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 17; this shifts rows 17:23 down to 18:24
# and keeps the dimension / formatting of the sheet consistent.
$ws.Rows.Item(17).Insert()

# Populate the newly inserted row 17 with the latest week's data.
$ws.Cells.Item(17, 1).Value  = 4
$ws.Cells.Item(17, 2).Value  = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(17, 3).Value  = "Los Lagos"
$ws.Cells.Item(17, 4).Value  = 44729
$ws.Cells.Item(17, 5).Value  = 10
$ws.Cells.Item(17, 6).Value  = 100112012
$ws.Cells.Item(17, 7).Value  = "Espinaca"
$ws.Cells.Item(17, 8).Value  = "Sin especificar"
$ws.Cells.Item(17, 9).Value  = "Primera"
$ws.Cells.Item(17, 10).Value = 35
$ws.Cells.Item(17, 11).Value = 13000
$ws.Cells.Item(17, 12).Value = 13000
$ws.Cells.Item(17, 13).Value = 13000
$ws.Cells.Item(17, 14).Value = "$/cuna 10 kilos"
$ws.Cells.Item(17, 15).Value = "Región Metropolitana"
$ws.Cells.Item(17, 16).Value = 1300
$ws.Cells.Item(17, 17).Value = 10
$ws.Cells.Item(17, 18).Value = "Hortaliza"
